$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each target cell is written with a leading escaped single-quote ("`'")
# so Excel treats the value as literal text (quote-prefix), preserving
# formatting such as leading zeros, repeated dots, and percent padding.
# The Style is reset to "Normal" afterward so no stray number-format /
# quote-prefix style index is left attached to the cell.

$ws.Range("D2").Value = "`'38.128.29"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "`'  +2.46%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "`'2.052.92"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "`'  +1.43%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "`'  +0.14%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "`'228.92"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "`'  +0.14%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("E6").Value = "`'  +0.86%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "`'60.59"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "`'  +7.78%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "`'  -0.06%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "`'0.385"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "`'  +1.81%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "`'0.0825"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "`'  +5.51%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "`'  +1.64%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "`'14.79"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "`'  +3.53%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "`'2.356.89"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "`'  +1.34%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "`'21.10"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "`'  +4.66%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "`'5.33"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "`'  +2.67%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "`'0.756"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "`'  +2.23%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "`'2.051.20"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "`'  +1.12%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "`'38.051.91"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "`'  +2.35%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "`'6.25"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "`'  +1.72%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "`'69.72"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "`'  +1.00%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "`'0.0₃0835"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "`'  +2.13%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "`'225.33"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "`'  +0.97%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "`'1.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "`'  +0.11%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "`'2.44"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "`'  -0.50%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "`'2.21"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "`'  -0.32%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "`'165.73"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "`'  +1.43%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "`'9.22"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "`'  +1.94%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "`'0.133"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "`'  +3.15%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "`'18.95"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "`'  +1.29%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = "`'  -1.44%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value = "`'  +2.44%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "`'4.51"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "`'  +1.52%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "`'4.57"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "`'  +2.65%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = "`'  +1.08%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "`'0.0604"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "`'  +0.20%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "`'6.41"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "`'  +16.50%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "`'2.28"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "`'  -2.38%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "`'3.26"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "`'  +1.79%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Value = "`'  +0.01%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "`'1.519.19"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "`'  +3.37%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("B41").Value = "`'InjectiveProtocol"
$ws.Range("B41").Style = "Normal"
$ws.Range("C41").Value = "`'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("C41").Style = "Normal"
$ws.Range("D41").Value = "`'16.96"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "`'  +4.60%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("B42").Value = "`'Aave"
$ws.Range("B42").Style = "Normal"
$ws.Range("C42").Value = "`'https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("C42").Style = "Normal"
$ws.Range("D42").Value = "`'97.12"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "`'  +3.33%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "`'  +0.74%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "`'2.86"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "`'  +1.88%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "`'0.0928"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "`'  +1.67%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "`'  +1.35%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "`'4.00"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "`'  -6.67%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "`'1.01"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "`'  +0.53%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "`'2.98"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "`'  +2.00%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "`'7.02"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "`'  -1.98%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "`'2.245.69"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "`'  +1.49%  "
$ws.Range("E51").Style = "Normal"
